$wb = $excel.ActiveWorkbook

# --- Sheet: zh-cn ---
$zh = $wb.Worksheets.Item("zh-cn")

# Status: "Handoff transform failed" -> "Ready for handoff"
$zh.Range("B2").Value = "Ready for handoff"

# Latest Handoff File: add hyperlink to the generated xlf handoff file
$zh.Hyperlinks.Add(
    $zh.Range("C2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/cef79d69a88e105989ac7e32f7ad405b8ae1fc93/e2e/7b2d39ff-d269-4e15-87bc-9b12eafbb1f9.5b5b2a7bf0de92d8d8f104cb87c11e077ba6e1f1.zh-cn.xlf",
    "",
    "",
    "7b2d39ff-d269-4e15-87bc-9b12eafbb1f9.5b5b2a7bf0de92d8d8f104cb87c11e077ba6e1f1.zh-cn.xlf"
)

# Latest Handoff Datetime
$zh.Range("D2").Value = "2016-01-27 08:38:35"

# Handoff Reason: "Ignored" -> "Include"
$zh.Range("H2").Value = "Include"

# --- Sheet: de-de ---
$de = $wb.Worksheets.Item("de-de")

# Status: "Handoff transform failed" -> "Ready for handoff"
$de.Range("B2").Value = "Ready for handoff"

# Latest Handoff File: add hyperlink to the generated xlf handoff file
$de.Hyperlinks.Add(
    $de.Range("C2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/cef79d69a88e105989ac7e32f7ad405b8ae1fc93/e2e/7b2d39ff-d269-4e15-87bc-9b12eafbb1f9.5b5b2a7bf0de92d8d8f104cb87c11e077ba6e1f1.de-de.xlf",
    "",
    "",
    "7b2d39ff-d269-4e15-87bc-9b12eafbb1f9.5b5b2a7bf0de92d8d8f104cb87c11e077ba6e1f1.de-de.xlf"
)

# Latest Handoff Datetime
$de.Range("D2").Value = "2016-01-27 08:38:46"

# Handoff Reason: "Ignored" -> "Include"
$de.Range("H2").Value = "Include"
